# Correção no modulo de geração de bombas
$wb = $excel.ActiveWorkbook

# --- Sheet "configurações": shrink board size from 5 to 3 ---
$wsConfig = $wb.Worksheets.Item("configurações")
$wsConfig.Range("B1").NumberFormat = "@"
$wsConfig.Range("B1").Value = "3"
$wsConfig.Range("B1").ClearFormats()
$wsConfig.Range("B2").NumberFormat = "@"
$wsConfig.Range("B2").Value = "3"
$wsConfig.Range("B2").ClearFormats()

# --- Sheet "jogo": resize bomb grid from 5x5 to 3x3 with new layout ---
$wsJogo = $wb.Worksheets.Item("jogo")

# Drop the now-unused columns D:E and rows 4:5 so the board becomes 3x3
$wsJogo.Range("D1:E5").EntireColumn.Delete()
$wsJogo.Range("A4:C5").EntireRow.Delete()

# Update the bomb layout for the remaining 3x3 grid
$wsJogo.Range("A2").Value = -1
$wsJogo.Range("C2").Value = 0
